# Fruta / hortaliza, semanal
# Shuffle the weekly data rows (2-10) of the Espárragos sheet so that each
# row now carries the values that belong to a different date, matching the
# updated weekly extract. Only columns D, H, J, K, L, M, N, O and P actually
# change value (the remaining columns are identical for every row already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (target) values for rows 2..10, columns D,H,J,K,L,M,N,O,P
$rows = @(
  @{ Row = 2;  D = 44526; H = "Sin especificar"; J = 100; K = 1500; L = 1600; M = 1550; N = "`$/kilo";    O = "Provincia de Linares"; P = 1550 },
  @{ Row = 3;  D = 44468; H = "Verde";           J = 500; K = 1800; L = 2000; M = 1920; N = "`$/kilo";    O = "Provincia de Linares"; P = 1920 },
  @{ Row = 4;  D = 44524; H = "Sin especificar"; J = 200; K = 1500; L = 1600; M = 1550; N = "`$/kilo";    O = "Provincia de Talca";   P = 1550 },
  @{ Row = 5;  D = 44511; H = "Sin especificar"; J = 600; K = 1300; L = 1400; M = 1350; N = "`$/kilo";    O = "Provincia de Linares"; P = 1350 },
  @{ Row = 6;  D = 44477; H = "Sin especificar"; J = 500; K = 1400; L = 1500; M = 1460; N = "`$/kilo";    O = "Provincia de Linares"; P = 1460 },
  @{ Row = 7;  D = 44510; H = "Sin especificar"; J = 600; K = 1300; L = 1400; M = 1350; N = "`$/kilo";    O = "Provincia de Linares"; P = 1350 },
  @{ Row = 8;  D = 44496; H = "Sin especificar"; J = 550; K = 1500; L = 2000; M = 1773; N = "`$/paquete"; O = "Provincia de Linares"; P = 1773 },
  @{ Row = 9;  D = 44519; H = "Sin especificar"; J = 250; K = 1200; L = 1300; M = 1240; N = "`$/kilo";    O = "Provincia de Linares"; P = 1240 },
  @{ Row = 10; D = 44489; H = "Sin especificar"; J = 600; K = 1400; L = 1500; M = 1450; N = "`$/kilo";    O = "Provincia de Linares"; P = 1450 }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
  $ws.Cells.Item($row, 8).Value  = $r.H   # H: Variedad
  $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
  $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio minimo
  $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio maximo
  $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
  $ws.Cells.Item($row, 14).Value = $r.N   # N: Unidad de comercializacion
  $ws.Cells.Item($row, 15).Value = $r.O   # O: Origen
  $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
